$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate row 2 into row 3 (preserves formatting/styles) before we touch values.
$ws.Range("A2:T2").Copy($ws.Range("A3:T3"))

# 2. Insert a new column before K ("Wins") for "Precision Crossing".
#    Existing K..T content (both rows) shifts right to L..U.
$ws.Columns("K").Insert()

# 3. Header for the new column.
$ws.Range("K1").Value = "Precision Crossing"

# 4. Update the "From" date for both test rows. Force text format first so
#    Excel doesn't auto-coerce the ISO-looking string into a date serial,
#    then drop back to the unstyled "Normal" style (matches original file,
#    which kept these as plain shared-string cells with no cell style).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2021-07-01"
$ws.Range("C2").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2021-07-01"
$ws.Range("C3").Style = "Normal"

# 5. Row 2 (Test #1): precision crossing OFF.
$ws.Range("K2").Value = $false
$ws.Range("L2").Value = 43
$ws.Range("M2").Value = 44
$ws.Range("N2").Value = 87
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "49.4%"
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").Value = -6
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 6450
$ws.Range("S2").Value = -4400
$ws.Range("T2").Value = 875.7125000000002
$ws.Range("U2").Value = 1174.2875

# 6. Row 3 (Test #2): precision crossing ON.
$ws.Range("A3").Value = 2
$ws.Range("K3").Value = $true
$ws.Range("L3").Value = 41
$ws.Range("M3").Value = 43
$ws.Range("N3").Value = 84
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "48.8%"
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = -5
$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 6150
$ws.Range("S3").Value = -4300
$ws.Range("T3").Value = 850.7125000000003
$ws.Range("U3").Value = 999.2874999999997
